$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for B:E and G columns, rows 2-8 (F column is unchanged)
$data = @{
    2 = @{ B = 1.455362044514542;  C = 0.04071648406533734; D = 0.1494219747398047; E = 0.4942365360607697; G = 2.139737039380454 }
    3 = @{ B = 0.6606524410359556; C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; G = 2.960089034096801 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    5 = @{ B = 0.6606524410359556; C = 0.002571899574220771; D = 0.7527432677738641; E = 10.19245300693656;  G = 11.6084206153206 }
    6 = @{ B = 0.6606524410359556; C = 1.655778082260271;    D = 3.537761648806719;  E = 10.19245300693656;  G = 16.0466451790395 }
    7 = @{ B = 0.2917716402565462; C = 1.655778082260271;    D = 22.3905356188092;   E = 10.19245300693656;  G = 34.53053834826257 }
    8 = @{ B = 0.1190320826869504; C = 0.306821227259698;    D = 0.7527432677738641; E = 0.4942365360607697; G = 1.672833113781282 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
